# The post "「象の群れ。ドローンの音から逃げている」" (row 775) was removed from the
# source data, so delete that entire row. Excel automatically shifts all
# subsequent rows (previously 776:852) up by one, producing the new
# range A1:C851.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(775).Delete()
